$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure the cells keep their original text (string) representation
# rather than being auto-converted to numbers/percentages by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "306.77"
$ws.Range("E2").Value = "0.75%"
$ws.Range("D3").Value = "36.22"
$ws.Range("E3").Value = "1.06%"
$ws.Range("D4").Value = "5.054"
$ws.Range("E4").Value = "-0.40%"
$ws.Range("D5").Value = "0.08063"
$ws.Range("E5").Value = "0.02%"
$ws.Range("D6").Value = "2.131"
$ws.Range("E6").Value = "10.16%"
$ws.Range("D7").Value = "7.832"
$ws.Range("E7").Value = "-0.10%"
$ws.Range("D8").Value = "0.9255"
$ws.Range("E8").Value = "-0.49%"
$ws.Range("D9").Value = "0.1446"
$ws.Range("E9").Value = "12.39%"
$ws.Range("D10").Value = "0.1916"
$ws.Range("E10").Value = "0.54%"
$ws.Range("D11").Value = "0.09078"
$ws.Range("E11").Value = "-1.18%"
$ws.Range("D12").Value = "0.03442"
$ws.Range("E12").Value = "-1.14%"
$ws.Range("D13").Value = "0.09908"
$ws.Range("E13").Value = "-0.02%"
$ws.Range("D14").Value = "0.001411"
$ws.Range("E14").Value = "-0.89%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "0.04358"
$ws.Range("E15").Value = "-1.21%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006098"
$ws.Range("E16").Value = "-8.60%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.841"
$ws.Range("E17").Value = "6.34%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "4.139"
$ws.Range("E18").Value = "-0.26%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "3.397"
$ws.Range("E19").Value = "11.77%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3449"
$ws.Range("E20").Value = "0.75%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "0.1335"
$ws.Range("E21").Value = "-0.13%"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "4.794"
$ws.Range("E22").Value = "-7.05%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "0.2338"
$ws.Range("E23").Value = "-7.71%"
$ws.Range("D24").Value = "0.001229"
$ws.Range("E24").Value = "-0.41%"
$ws.Range("D25").Value = "0.004307"
$ws.Range("E25").Value = "-8.27%"
$ws.Range("D27").Value = "0.0001300"
$ws.Range("E27").Value = "-0.01%"
$ws.Range("D39").Value = "0.02021"
$ws.Range("E39").Value = "1.43%"
$ws.Range("D40").Value = "0.05140"
$ws.Range("E40").Value = "-1.04%"
$ws.Range("D41").Value = "0.007517"
$ws.Range("E41").Value = "-1.34%"
$ws.Range("D42").Value = "0.01011"
$ws.Range("E42").Value = "-0.07%"
$ws.Range("E43").Value = "-0.22%"
$ws.Range("D44").Value = "0.002151"
$ws.Range("E44").Value = "2.37%"
$ws.Range("D45").Value = "0.009958"
$ws.Range("E45").Value = "-6.97%"
$ws.Range("D46").Value = "0.00006274"
$ws.Range("E46").Value = "-0.38%"
$ws.Range("E47").Value = "-0.17%"
$ws.Range("D48").Value = "63.75"
$ws.Range("E48").Value = "0.29%"
$ws.Range("D49").Value = "0.001248"
$ws.Range("E49").Value = "-22.08%"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("E50").Value = "-0.17%"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").Value = "-0.17%"
